# Weekly update: insert a new week's row of data at the top of the
# Perejil (parsley) price series (row 216) and shift the existing
# history down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 216, pushing rows 216:255 down to 217:256.
# xlShiftDown = -4121
$ws.Rows("216:216").Insert(-4121)

# The row that used to be 216 is now 217; re-use its (constant) values
# for the new row, then overwrite the date (D) and volume (J) with the
# new week's figures.
$ws.Range("A216").Value = $ws.Range("A217").Value()
$ws.Range("B216").Value = $ws.Range("B217").Value()
$ws.Range("C216").Value = $ws.Range("C217").Value()
$ws.Range("D216").Value = 45209
$ws.Range("E216").Value = $ws.Range("E217").Value()
$ws.Range("F216").Value = $ws.Range("F217").Value()
$ws.Range("G216").Value = $ws.Range("G217").Value()
$ws.Range("H216").Value = $ws.Range("H217").Value()
$ws.Range("I216").Value = $ws.Range("I217").Value()
$ws.Range("J216").Value = 2600
$ws.Range("K216").Value = $ws.Range("K217").Value()
$ws.Range("L216").Value = $ws.Range("L217").Value()
$ws.Range("M216").Value = $ws.Range("M217").Value()
$ws.Range("N216").Value = $ws.Range("N217").Value()
$ws.Range("O216").Value = $ws.Range("O217").Value()
$ws.Range("P216").Value = $ws.Range("P217").Value()
$ws.Range("Q216").Value = $ws.Range("Q217").Value()
$ws.Range("R216").Value = $ws.Range("R217").Value()

# Keep the date column's number format consistent with the rest of
# the series (YYYY-MM-DD HH:MM:SS).
$ws.Range("D216").NumberFormat = $ws.Range("D217").NumberFormat()
